# FormatoPF.xlsx - "Codigo corregido de detalles pequeños"
#
# 1) Two shared-text typos are shortened ("dos" -> "do", "tres" -> "tr")
#    on the "Origen Razón Social Domicilio" column (AZ) for the 2nd and
#    3rd data rows.
# 2) "Numero dependientes" (U3) is corrected from 1 to 12.
# 3) A new narrow spacer column is introduced right after the existing
#    data block (column 100 / CV) and is stamped with a single blank
#    space on each of the four detail rows (3-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1) Fix the truncated words -------------------------------------------
$ws.Range("AZ3").Value = "do"
$ws.Range("AZ4").Value = "tr"

# -- 2) Correct the dependents count on row 3 ------------------------------
$ws.Range("U3").Value = 12

# -- 3) New narrow spacer column (100 / CV) --------------------------------
$ws.Columns.Item(100).ColumnWidth = 0.6

$ws.Range("CV3").Value = " "
$ws.Range("CV4").Value = " "
$ws.Range("CV5").Value = " "
$ws.Range("CV6").Value = " "
